$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6372756958007812
$ws.Range("B1").Value = 1.485128164291382
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.295554637908936
$ws.Range("E1").Value = 1.365162253379822
